$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reward")

# --- Set cell values ---
$ws.Cells.Item(1,1).Value = "ID"
$ws.Cells.Item(1,2).Value = "Desc"
$ws.Cells.Item(1,3).Value = "RewardItem1Id"
$ws.Cells.Item(1,4).Value = "RewardItem1Num"
$ws.Cells.Item(1,5).Value = "RewardItem2Id"
$ws.Cells.Item(1,6).Value = "RewardItem2Num"
$ws.Cells.Item(1,7).Value = "ExchangeId"
$ws.Cells.Item(1,8).Value = "InputItem1Id"
$ws.Cells.Item(1,9).Value = "InputItem1Num"
$ws.Cells.Item(1,10).Value = "InputItem2Id"
$ws.Cells.Item(1,11).Value = "InputItem2Num"
$ws.Cells.Item(1,12).Value = "InputItem3Id"
$ws.Cells.Item(1,13).Value = "InputItem3Num"
$ws.Cells.Item(1,14).Value = "OutputItem1Id"
$ws.Cells.Item(1,15).Value = "OutputItem1Num"
$ws.Cells.Item(1,16).Value = "OutputItem2Id"
$ws.Cells.Item(1,17).Value = "OutputItem2Num"
$ws.Cells.Item(2,1).Value = "map<uint32, Reward>"
$ws.Cells.Item(2,2).Value = "string"
$ws.Cells.Item(2,3).Value = "[Item]int32"
$ws.Cells.Item(2,4).Value = "int32"
$ws.Cells.Item(2,5).Value = "int32"
$ws.Cells.Item(2,6).Value = "int32"
$ws.Cells.Item(2,7).Value = "[Exchange]uint32"
$ws.Cells.Item(2,8).Value = "[Item]int32"
$ws.Cells.Item(2,9).Value = "int32"
$ws.Cells.Item(2,10).Value = "int32"
$ws.Cells.Item(2,11).Value = "int32"
$ws.Cells.Item(2,12).Value = "int32"
$ws.Cells.Item(2,13).Value = "int32"
$ws.Cells.Item(2,14).Value = "[Item]int32"
$ws.Cells.Item(2,15).Value = "int32"
$ws.Cells.Item(2,16).Value = "int32"
$ws.Cells.Item(2,17).Value = "int32"
$ws.Cells.Item(3,1).Value = "奖励ID"
$ws.Cells.Item(3,2).Value = "描述"
$ws.Cells.Item(3,3).Value = "奖励1Id"
$ws.Cells.Item(3,4).Value = "奖励1Num"
$ws.Cells.Item(3,5).Value = "奖励2Id"
$ws.Cells.Item(3,6).Value = "奖励2Num"
$ws.Cells.Item(3,7).Value = "兑换ID"
$ws.Cells.Item(3,8).Value = "道具1Id"
$ws.Cells.Item(3,9).Value = "道具1Num"
$ws.Cells.Item(3,10).Value = "道具2Id"
$ws.Cells.Item(3,11).Value = "道具2Num"
$ws.Cells.Item(3,12).Value = "道具3Id"
$ws.Cells.Item(3,13).Value = "道具3Num"
$ws.Cells.Item(3,14).Value = "奖励1Id"
$ws.Cells.Item(3,15).Value = "奖励1Num"
$ws.Cells.Item(3,16).Value = "奖励2Id"
$ws.Cells.Item(3,17).Value = "奖励2Num"
$ws.Cells.Item(4,1).Value = 1
$ws.Cells.Item(4,2).Value = "奖励1"
$ws.Cells.Item(4,3).Value = 1
$ws.Cells.Item(4,4).Value = 2002
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 2002
$ws.Cells.Item(4,7).Value = 1
$ws.Cells.Item(4,8).Value = 2002
$ws.Cells.Item(4,9).Value = 3
$ws.Cells.Item(4,10).Value = 2002
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 2002
$ws.Cells.Item(4,13).Value = 3
$ws.Cells.Item(4,14).Value = 2002
$ws.Cells.Item(4,15).Value = 3
$ws.Cells.Item(5,1).Value = 1
$ws.Cells.Item(5,2).Value = "奖励2"
$ws.Cells.Item(5,3).Value = 2
$ws.Cells.Item(5,4).Value = 2002
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 2002
$ws.Cells.Item(5,7).Value = 2
$ws.Cells.Item(5,8).Value = 2002
$ws.Cells.Item(5,9).Value = 3
$ws.Cells.Item(5,10).Value = 2002
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 2002
$ws.Cells.Item(5,13).Value = 3
$ws.Cells.Item(5,14).Value = 2002
$ws.Cells.Item(5,15).Value = 3
$ws.Cells.Item(6,1).Value = 1
$ws.Cells.Item(6,2).Value = "奖励3"
$ws.Cells.Item(6,3).Value = 3
$ws.Cells.Item(6,4).Value = 2002
$ws.Cells.Item(6,7).Value = 2
$ws.Cells.Item(6,8).Value = 2002
$ws.Cells.Item(6,9).Value = 3
$ws.Cells.Item(6,10).Value = 2002
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 2002
$ws.Cells.Item(6,13).Value = 3
$ws.Cells.Item(6,14).Value = 2002
$ws.Cells.Item(6,15).Value = 3
$ws.Cells.Item(7,1).Value = 2
$ws.Cells.Item(7,2).Value = "奖励3"
$ws.Cells.Item(7,3).Value = 1
$ws.Cells.Item(7,4).Value = 2001
$ws.Cells.Item(7,7).Value = 1
$ws.Cells.Item(7,8).Value = 2007
$ws.Cells.Item(7,9).Value = 10
$ws.Cells.Item(7,10).Value = 2001
$ws.Cells.Item(7,11).Value = 1
$ws.Cells.Item(7,14).Value = 2007
$ws.Cells.Item(7,15).Value = 10

# --- Copy formats for newly extended columns/cells ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("I1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("I2:Q2").PasteSpecial(-4122) | Out-Null
$ws.Range("A3").Copy() | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("H3:Q3").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# --- Column widths (approximate bestFit target widths) ---
$ws.Columns.Item(3).ColumnWidth = 13.25
$ws.Columns.Item(4).ColumnWidth = 15.5
$ws.Columns.Item(5).ColumnWidth = 13.25
$ws.Columns.Item(6).ColumnWidth = 15.5
$ws.Columns.Item(7).ColumnWidth = 15.125
$ws.Columns.Item(8).ColumnWidth = 11.125
$ws.Columns.Item(9).ColumnWidth = 13.5
$ws.Columns.Item(10).ColumnWidth = 11.125
$ws.Columns.Item(11).ColumnWidth = 13.5
$ws.Columns.Item(12).ColumnWidth = 11.125
$ws.Columns.Item(13).ColumnWidth = 13.5
$ws.Columns.Item(14).ColumnWidth = 12.625
$ws.Columns.Item(15).ColumnWidth = 15
$ws.Columns.Item(16).ColumnWidth = 12.625
$ws.Columns.Item(17).ColumnWidth = 15

# --- Selection ---
$ws.Range("H14").Select() | Out-Null
